$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New registrations appended below the existing rows (rows 2-5), growing
# the used range from A1:D5 to A1:D9.
$data = @(
    @("myemail342@gmail.com",      "oi",     "2024-12-01", "12:21:59"),
    @("myemailyfguh@gmail.com",    "12",     "2024-12-01", "15:20:04"),
    @("criando_conta@gmail.com",   "12",     "2024-12-01", "15:20:48"),
    @("criando_conta123@gmail.com","123123", "2024-12-01", "15:21:56")
)

$row = 6
foreach ($r in $data) {
    $emailCell = $ws.Cells.Item($row, 1)
    $pwdCell   = $ws.Cells.Item($row, 2)
    $dateCell  = $ws.Cells.Item($row, 3)
    $timeCell  = $ws.Cells.Item($row, 4)

    # Password and date columns can look like pure numbers/dates
    # ("12", "123123", "2024-12-01") -- force them to plain text so Excel
    # doesn't silently convert them to a Number/Date cell, matching the
    # original inline-string text cells used throughout this sheet. The
    # text number format is removed again right after so the cell keeps
    # its default (unstyled) formatting, same as every other cell here.
    $pwdCell.NumberFormat = "@"
    $dateCell.NumberFormat = "@"

    $emailCell.Value = $r[0]
    $pwdCell.Value = $r[1]
    $dateCell.Value = $r[2]
    $timeCell.Value = $r[3]

    $pwdCell.ClearFormats()
    $dateCell.ClearFormats()

    $row++
}
